{"js": "// Update the two-digit-division answer table: replace each old\n// \"a\u00f7b=c, d\" answer string with its corresponding new value.\n// Every cell text in the document is unique, so body.search() for the\n// exact old text unambiguously locates the run to replace.\n\nconst replacements = [\n  [\"32\u00f78=4, 0\", \"93\u00f73=31, 0\"],\n  [\"55\u00f74=13, 3\", \"65\u00f77=9, 2\"],\n  [\"16\u00f77=2, 2\", \"15\u00f72=7, 1\"],\n  [\"17\u00f75=3, 2\", \"49\u00f76=8, 1\"],\n  [\"44\u00f79=4, 8\", \"93\u00f76=15, 3\"],\n  [\"73\u00f74=18, 1\", \"49\u00f74=12, 1\"],\n  [\"40\u00f73=13, 1\", \"65\u00f77=9, 2\"],\n  [\"37\u00f72=18, 1\", \"92\u00f77=13, 1\"],\n  [\"44\u00f73=14, 2\", \"79\u00f72=39, 1\"],\n  [\"22\u00f74=5, 2\", \"27\u00f79=3, 0\"],\n  [\"56\u00f74=14, 0\", \"65\u00f76=10, 5\"],\n  [\"77\u00f75=15, 2\", \"46\u00f76=7, 4\"],\n  [\"57\u00f76=9, 3\", \"43\u00f79=4, 7\"],\n  [\"97\u00f79=10, 7\", \"35\u00f73=11, 2\"],\n  [\"59\u00f74=14, 3\", \"30\u00f76=5, 0\"],\n  [\"90\u00f75=18, 0\", \"29\u00f78=3, 5\"],\n  [\"61\u00f72=30, 1\", \"92\u00f78=11, 4\"],\n  [\"10\u00f78=1, 2\", \"57\u00f77=8, 1\"],\n  [\"74\u00f77=10, 4\", \"94\u00f79=10, 4\"],\n  [\"38\u00f74=9, 2\", \"58\u00f78=7, 2\"],\n  [\"33\u00f77=4, 5\", \"57\u00f75=11, 2\"],\n  [\"14\u00f78=1, 6\", \"69\u00f78=8, 5\"],\n  [\"96\u00f73=32, 0\", \"54\u00f74=13, 2\"],\n  [\"73\u00f78=9, 1\", \"96\u00f77=13, 5\"],\n  [\"16\u00f73=5, 1\", \"72\u00f73=24, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-division answer table: replace each old\n# \"a\u00f7b=c, d\" answer string with its corresponding new value using\n# Find/Replace over the document content. Every cell text in the\n# document is unique, so an exact, case-sensitive, non-wildcard find\n# unambiguously locates the text to replace.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"32\u00f78=4, 0\", \"93\u00f73=31, 0\"),\n  @(\"55\u00f74=13, 3\", \"65\u00f77=9, 2\"),\n  @(\"16\u00f77=2, 2\", \"15\u00f72=7, 1\"),\n  @(\"17\u00f75=3, 2\", \"49\u00f76=8, 1\"),\n  @(\"44\u00f79=4, 8\", \"93\u00f76=15, 3\"),\n  @(\"73\u00f74=18, 1\", \"49\u00f74=12, 1\"),\n  @(\"40\u00f73=13, 1\", \"65\u00f77=9, 2\"),\n  @(\"37\u00f72=18, 1\", \"92\u00f77=13, 1\"),\n  @(\"44\u00f73=14, 2\", \"79\u00f72=39, 1\"),\n  @(\"22\u00f74=5, 2\", \"27\u00f79=3, 0\"),\n  @(\"56\u00f74=14, 0\", \"65\u00f76=10, 5\"),\n  @(\"77\u00f75=15, 2\", \"46\u00f76=7, 4\"),\n  @(\"57\u00f76=9, 3\", \"43\u00f79=4, 7\"),\n  @(\"97\u00f79=10, 7\", \"35\u00f73=11, 2\"),\n  @(\"59\u00f74=14, 3\", \"30\u00f76=5, 0\"),\n  @(\"90\u00f75=18, 0\", \"29\u00f78=3, 5\"),\n  @(\"61\u00f72=30, 1\", \"92\u00f78=11, 4\"),\n  @(\"10\u00f78=1, 2\", \"57\u00f77=8, 1\"),\n  @(\"74\u00f77=10, 4\", \"94\u00f79=10, 4\"),\n  @(\"38\u00f74=9, 2\", \"58\u00f78=7, 2\"),\n  @(\"33\u00f77=4, 5\", \"57\u00f75=11, 2\"),\n  @(\"14\u00f78=1, 6\", \"69\u00f78=8, 5\"),\n  @(\"96\u00f73=32, 0\", \"54\u00f74=13, 2\"),\n  @(\"73\u00f78=9, 1\", \"96\u00f77=13, 5\"),\n  @(\"16\u00f73=5, 1\", \"72\u00f73=24, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
